$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values per diff
$ws.Range("D12").Value = 13092
$ws.Range("C13").Value = 407
$ws.Range("D13").Value = 13600

# Add new row 14 (force B14 to remain text, not auto-converted to a date)
$ws.Range("A14").Value = "Sergipe"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "01/01/2022"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = 393
$ws.Range("D14").Value = 15447
